# Fruta / hortaliza, semanal
#
# Two new daily price records are inserted above the existing data block
# (the data table starts at row 2; the block being extended runs through
# the former last row, 310). The two new rows push all the existing
# records for rows 235-310 down by two rows (to 237-312), and the sheet's
# used range grows from A1:R310 to A1:R312.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 235, shifting rows 235:310 down to 237:312.
$ws.Range("A235:A236").EntireRow.Insert()

# New row 235: Femacal de La Calera, Provincia de Quillota, 36-unit box.
$ws.Cells.Item(235, 1).Value = 3
$ws.Cells.Item(235, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(235, 3).Value = "Coquimbo"
$ws.Cells.Item(235, 4).Value = 44559
$ws.Cells.Item(235, 5).Value = 5
$ws.Cells.Item(235, 6).Value = 100112032
$ws.Cells.Item(235, 7).Value = "Zapallo italiano"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 135
$ws.Cells.Item(235, 11).Value = 4000
$ws.Cells.Item(235, 12).Value = 4500
$ws.Cells.Item(235, 13).Value = 4252
$ws.Cells.Item(235, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(235, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(235, 16).Value = 118
$ws.Cells.Item(235, 17).Value = 36
$ws.Cells.Item(235, 18).Value = "Hortaliza"

# New row 236: Femacal de La Calera, Provincia de Quillota, 70-unit box.
$ws.Cells.Item(236, 1).Value = 3
$ws.Cells.Item(236, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44559
$ws.Cells.Item(236, 5).Value = 5
$ws.Cells.Item(236, 6).Value = 100112032
$ws.Cells.Item(236, 7).Value = "Zapallo italiano"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 148
$ws.Cells.Item(236, 11).Value = 6500
$ws.Cells.Item(236, 12).Value = 7000
$ws.Cells.Item(236, 13).Value = 6736
$ws.Cells.Item(236, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(236, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(236, 16).Value = 96
$ws.Cells.Item(236, 17).Value = 70
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# Keep the date columns formatted like the rest of the table (style index 2
# carries the date/time number format used throughout column D).
$dateFormat = $ws.Cells.Item(234, 4).NumberFormat
$ws.Cells.Item(235, 4).NumberFormat = $dateFormat
$ws.Cells.Item(236, 4).NumberFormat = $dateFormat
